# Auto-generated edit script: update Leve profit market-price figures across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 11171
$ws.Range("I40").Value = 21676.2
$ws.Range("J40").Value = 2416.6667
$ws.Range("K40").Value = 21676.2
$ws.Range("L40").Value = 2416.6667
$ws.Range("M40").Value = -21501.2
$ws.Range("N40").Value = -2766.6667
$ws.Range("H51").Value = 8550300
$ws.Range("I51").Value = 12349556
$ws.Range("J51").Value = 1975
$ws.Range("K51").Value = 12349556
$ws.Range("L51").Value = 1975
$ws.Range("M51").Value = -12349072
$ws.Range("N51").Value = -2943
$ws.Range("H69").Value = 3249.5
$ws.Range("H72").Value = 3249.5
$ws.Range("H112").Value = 1954.7693
$ws.Range("I112").Value = 1900
$ws.Range("J112").Value = 1961.9131
$ws.Range("K112").Value = 5700
$ws.Range("L112").Value = 5885.7393
$ws.Range("M112").Value = -4592
$ws.Range("N112").Value = -8101.7393
$ws.Range("H114").Value = 45450
$ws.Range("J114").Value = 45450
$ws.Range("L114").Value = 45450
$ws.Range("N114").Value = -54128
$ws.Range("H129").Value = 2297.625
$ws.Range("J129").Value = 3314.5833
$ws.Range("L129").Value = 9943.749899999999
$ws.Range("N129").Value = -19943.7499
$ws.Range("H141").Value = 2100.1428
$ws.Range("I141").Value = 811.69696
$ws.Range("J141").Value = 6824.4443
$ws.Range("K141").Value = 2435.09088
$ws.Range("L141").Value = 20473.3329
$ws.Range("M141").Value = 2744.90912
$ws.Range("N141").Value = -30833.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2397.5625
$ws.Range("I2").Value = 2497.2144
$ws.Range("J2").Value = 1700
$ws.Range("K2").Value = 2497.2144
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = -2384.2144
$ws.Range("N2").Value = -1926
$ws.Range("H61").Value = 1602.9512
$ws.Range("I61").Value = 1367.3889
$ws.Range("K61").Value = 1367.3889
$ws.Range("M61").Value = -1155.3889
$ws.Range("H74").Value = 1452.6364
$ws.Range("I74").Value = 1331.36
$ws.Range("J74").Value = 1831.625
$ws.Range("K74").Value = 1331.36
$ws.Range("L74").Value = 1831.625
$ws.Range("M74").Value = -457.3599999999999
$ws.Range("N74").Value = -3579.625
$ws.Range("H77").Value = 1452.6364
$ws.Range("I77").Value = 1331.36
$ws.Range("J77").Value = 1831.625
$ws.Range("K77").Value = 6656.799999999999
$ws.Range("L77").Value = 9158.125
$ws.Range("M77").Value = -2288.799999999999
$ws.Range("N77").Value = -17894.125
$ws.Range("H110").Value = 1495.6129
$ws.Range("I110").Value = 1498.56
$ws.Range("J110").Value = 1483.3334
$ws.Range("K110").Value = 1498.56
$ws.Range("L110").Value = 1483.3334
$ws.Range("M110").Value = 546.4400000000001
$ws.Range("N110").Value = -5573.3334
$ws.Range("H116").Value = 2397.5625
$ws.Range("I116").Value = 2497.2144
$ws.Range("J116").Value = 1700
$ws.Range("K116").Value = 2497.2144
$ws.Range("L116").Value = 1700
$ws.Range("M116").Value = -203.2143999999998
$ws.Range("N116").Value = -6288
$ws.Range("H122").Value = 1575.9056
$ws.Range("I122").Value = 1550.46
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4651.38
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2201.38
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 8476101
$ws.Range("I132").Value = 10870487
$ws.Range("J132").Value = 3655
$ws.Range("K132").Value = 32611461
$ws.Range("L132").Value = 10965
$ws.Range("M132").Value = -32608931
$ws.Range("N132").Value = -16025
$ws.Range("H134").Value = 45211.11
$ws.Range("J134").Value = 45211.11
$ws.Range("L134").Value = 45211.11
$ws.Range("N134").Value = -55351.11
$ws.Range("H136").Value = 1602.9512
$ws.Range("I136").Value = 1367.3889
$ws.Range("K136").Value = 4102.1667
$ws.Range("M136").Value = -1552.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2397.5625
$ws.Range("I3").Value = 2497.2144
$ws.Range("J3").Value = 1700
$ws.Range("K3").Value = 2497.2144
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = -2383.2144
$ws.Range("N3").Value = -1928
$ws.Range("H105").Value = 2183.2856
$ws.Range("I105").Value = 2465
$ws.Range("J105").Value = 2136.3333
$ws.Range("K105").Value = 2465
$ws.Range("L105").Value = 2136.3333
$ws.Range("M105").Value = -718
$ws.Range("N105").Value = -5630.3333
$ws.Range("H107").Value = 2144.9355
$ws.Range("I107").Value = 2058.6667
$ws.Range("K107").Value = 2058.6667
$ws.Range("M107").Value = -138.6667000000002
$ws.Range("H124").Value = 52328
$ws.Range("J124").Value = 52328
$ws.Range("L124").Value = 52328
$ws.Range("N124").Value = -62148
$ws.Range("H134").Value = 3190.307
$ws.Range("I134").Value = 1436.3667
$ws.Range("J134").Value = 4097.517
$ws.Range("K134").Value = 4309.1001
$ws.Range("L134").Value = 12292.551
$ws.Range("M134").Value = -1774.1001
$ws.Range("N134").Value = -17362.551

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 930.7692
$ws.Range("I16").Value = 922
$ws.Range("J16").Value = 960
$ws.Range("K16").Value = 922
$ws.Range("L16").Value = 960
$ws.Range("M16").Value = -635
$ws.Range("N16").Value = -1534
$ws.Range("H31").Value = 2223.38
$ws.Range("I31").Value = 1006.4889
$ws.Range("K31").Value = 1006.4889
$ws.Range("M31").Value = -711.4888999999999
$ws.Range("H34").Value = 2223.38
$ws.Range("I34").Value = 1006.4889
$ws.Range("K34").Value = 1006.4889
$ws.Range("M34").Value = -804.4888999999999
$ws.Range("H88").Value = 37998.375
$ws.Range("J88").Value = 37998.375
$ws.Range("L88").Value = 37998.375
$ws.Range("N88").Value = -38810.375
$ws.Range("H91").Value = 37998.375
$ws.Range("J91").Value = 37998.375
$ws.Range("L91").Value = 37998.375
$ws.Range("N91").Value = -40806.375
$ws.Range("H113").Value = 930.7692
$ws.Range("I113").Value = 922
$ws.Range("J113").Value = 960
$ws.Range("K113").Value = 922
$ws.Range("L113").Value = 960
$ws.Range("M113").Value = 1248
$ws.Range("N113").Value = -5300
$ws.Range("H132").Value = 50641.69
$ws.Range("I132").Value = 2139.7827
$ws.Range("J132").Value = 236565.67
$ws.Range("K132").Value = 6419.348100000001
$ws.Range("L132").Value = 709697.01
$ws.Range("M132").Value = -3889.348100000001
$ws.Range("N132").Value = -714757.01
$ws.Range("H134").Value = 304948.06
$ws.Range("I134").Value = 1163.4762
$ws.Range("J134").Value = 836571.0600000001
$ws.Range("K134").Value = 3490.4286
$ws.Range("L134").Value = 2509713.18
$ws.Range("M134").Value = -955.4286000000002
$ws.Range("N134").Value = -2514783.18

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4300
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 4300
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H114").Value = 1677.1
$ws.Range("I114").Value = 1006.8333
$ws.Range("J114").Value = 2682.5
$ws.Range("K114").Value = 3020.4999
$ws.Range("L114").Value = 8047.5
$ws.Range("M114").Value = 233.5001000000002
$ws.Range("N114").Value = -14555.5
$ws.Range("H131").Value = 3613.238
$ws.Range("J131").Value = 1434.3823
$ws.Range("L131").Value = 4303.1469
$ws.Range("N131").Value = -14383.1469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 170012.53
$ws.Range("I80").Value = 266456.3
$ws.Range("J80").Value = 3427.818
$ws.Range("K80").Value = 266456.3
$ws.Range("L80").Value = 3427.818
$ws.Range("M80").Value = -265458.3
$ws.Range("N80").Value = -5423.818
$ws.Range("H83").Value = 170012.53
$ws.Range("I83").Value = 266456.3
$ws.Range("J83").Value = 3427.818
$ws.Range("K83").Value = 1332281.5
$ws.Range("L83").Value = 17139.09
$ws.Range("M83").Value = -1327289.5
$ws.Range("N83").Value = -27123.09
$ws.Range("H123").Value = 13889.875
$ws.Range("J123").Value = 13889.875
$ws.Range("L123").Value = 13889.875
$ws.Range("N123").Value = -18789.875
$ws.Range("H141").Value = 36199.2
$ws.Range("J141").Value = 36199.2
$ws.Range("L141").Value = 36199.2
$ws.Range("N141").Value = -46559.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2856.85
$ws.Range("I40").Value = 2384.4666
$ws.Range("J40").Value = 4274
$ws.Range("K40").Value = 2384.4666
$ws.Range("L40").Value = 4274
$ws.Range("M40").Value = -2248.4666
$ws.Range("N40").Value = -4546
$ws.Range("H132").Value = 1791.3086
$ws.Range("I132").Value = 1259.2188
$ws.Range("J132").Value = 3794.4707
$ws.Range("K132").Value = 3777.6564
$ws.Range("L132").Value = 11383.4121
$ws.Range("M132").Value = -1247.6564
$ws.Range("N132").Value = -16443.4121
$ws.Range("H136").Value = 1742.5428
$ws.Range("I136").Value = 1449.6538
$ws.Range("J136").Value = 2588.6667
$ws.Range("K136").Value = 4348.9614
$ws.Range("L136").Value = 7766.000100000001
$ws.Range("M136").Value = -1798.9614
$ws.Range("N136").Value = -12866.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1524.64
$ws.Range("I132").Value = 1484.8529
$ws.Range("J132").Value = 1609.1875
$ws.Range("K132").Value = 4454.5587
$ws.Range("L132").Value = 4827.5625
$ws.Range("M132").Value = -1924.5587
$ws.Range("N132").Value = -9887.5625
